$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Avverkningsanmälningar")

# Column C ("Förändrad") holds the same date serial value (45190) for every
# data row (rows 2 through 453). The commit updates all of these to 45192.
$ws.Range("C2:C453").Value = 45192
